$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Report date: 2025-12-03 -> 2025-12-05 (stored as plain text, not a date) ---
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2025-12-05"
$ws.Range("A2").Style = "Normal"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-12-05"
$ws.Range("A3").Style = "Normal"

# --- MACRO_SIGNAL text: "🟢 상승 우위 (다소 완화)" -> "⚪ 중립 구간" ---
$ws.Range("O2").Value = "⚪ 중립 구간"
$ws.Range("O3").Value = "⚪ 중립 구간"

# --- Row 2 (Oklo Inc. / OKLO) metric updates ---
$ws.Range("D2").Value = 109.53
$ws.Range("E2").Value = 56
$ws.Range("F2").Value = 23.46
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 73
$ws.Range("J2").Value = 86
$ws.Range("K2").Value = 58.7
$ws.Range("N2").Value = 52.43913937059539

# --- Row 3 (NuScale Power Corporation / SMR) metric updates ---
$ws.Range("D3").Value = 22.7
$ws.Range("E3").Value = 48.6
$ws.Range("F3").Value = 19.25
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 50
$ws.Range("J3").Value = 53
$ws.Range("K3").Value = 54.7
$ws.Range("N3").Value = 52.43913937059539
